# Update row 6 of the "Ds cấp phát TSCĐ" sheet with the new asset entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ds cấp phát TSCĐ")

$ws.Range("B6").Value = "17/05/2025"
$ws.Range("C6").Value = "TS-008196"
$ws.Range("D6").Value = "KHO_VTYT"
$ws.Range("E6").Value = "A1"
